$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the duplicated/typo'd duration label "One days" -> "One day" in the
# "Estimated time to complete" column for the Composite OLAP rows. Since
# "One day" already exists elsewhere in the sheet, this removes the now
# -unused duplicate shared string on save.
$ws.Range("E12").Value = "One day"
$ws.Range("E13").Value = "One day"
$ws.Range("E14").Value = "One day"
$ws.Range("E15").Value = "One day"

# Move the active selection to I18 (single cell), matching where the user
# left the cursor after making the edit.
$ws.Range("I18").Select()
